# Natmi following Dr Hou advice
# Re-run of the NATMI ligand-receptor analysis (Npy -> Npy1r) adds the
# "ECs" and "FAPs" clusters alongside the existing "M2" cluster, expanding
# the single result row into six sending/target cluster combinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    # Sending, Ligand, Receptor, Target, then the 16 numeric metric columns (E:T)
    @("ECs", "Npy", "Npy1r", "ECs",
      1, 0.3333333333333333, 0.672624, 2.017872, 0.5229110735646243, 0.5229110735646243,
      2, 0.6666666666666666, 1.150707666666667, 3.452123, 0.1434406019171491, 0.1434406019171491,
      0.773993593584, 6.965942342256, 0.07500667914125236, 0.07500667914125236),
    @("ECs", "Npy", "Npy1r", "FAPs",
      1, 0.3333333333333333, 0.672624, 2.017872, 0.5229110735646243, 0.5229110735646243,
      1, 0.3333333333333333, 0.0331, 0.0993, 0.00412605569684884, 0.00412605569684884,
      0.0222638544, 0.2003746896, 0.002157560214026661, 0.002157560214026661),
    @("ECs", "Npy", "Npy1r", "sCs",
      1, 0.3333333333333333, 0.672624, 2.017872, 0.5229110735646243, 0.5229110735646243,
      3, 1, 6.838381666666667, 20.515145, 0.852433342386002, 0.852433342386002,
      4.59965963016, 41.39693667144, 0.4457468342093452, 0.4457468342093452),
    @("M2", "Npy", "Npy1r", "ECs",
      3, 1, 0.6136826666666667, 1.841048, 0.4770889264353757, 0.4770889264353757,
      2, 0.6666666666666666, 1.150707666666667, 3.452123, 0.1434406019171491, 0.1434406019171491,
      0.7061693494337778, 6.355524144904001, 0.06843392277589677, 0.06843392277589677),
    @("M2", "Npy", "Npy1r", "FAPs",
      3, 1, 0.6136826666666667, 1.841048, 0.4770889264353757, 0.4770889264353757,
      1, 0.3333333333333333, 0.0331, 0.0993, 0.00412605569684884, 0.00412605569684884,
      0.02031289626666666, 0.1828160664, 0.001968495482822179, 0.001968495482822179),
    @("M2", "Npy", "Npy1r", "sCs",
      3, 1, 0.6136826666666667, 1.841048, 0.4770889264353757, 0.4770889264353757,
      3, 1, 6.838381666666667, 20.515145, 0.852433342386002, 0.852433342386002,
      4.196596296884445, 37.76936667196001, 0.4066865081766567, 0.4066865081766567)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $values = $rows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}
